$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1), columns B..E ----
$ws.Range("B1").Value = "SP.RUR.TOTL.ZS:VGB"
$ws.Range("C1").Value = "SP.RUR.TOTL:VGB"
$ws.Range("D1").Value = "SP.URB.TOTL.IN.ZS:VGB"
$ws.Range("E1").Value = "SP.URB.TOTL:VGB"

# ---- Row labels (column A), rows 2..9 ----
$ws.Range("A2").Value = "SP.POP.TOTL:VGB:cor-value"
$ws.Range("A3").Value = "SP.POP.TOTL:VGB:p-value"
$ws.Range("A4").Value = "SP.RUR.TOTL.ZS:VGB:cor-value"
$ws.Range("A5").Value = "SP.RUR.TOTL.ZS:VGB:p-value"
$ws.Range("A6").Value = "SP.RUR.TOTL:VGB:cor-value"
$ws.Range("A7").Value = "SP.RUR.TOTL:VGB:p-value"
$ws.Range("A8").Value = "SP.URB.TOTL:VGB:cor-value"
$ws.Range("A9").Value = "SP.URB.TOTL:VGB:p-value"

# ---- Data values ----
# Scientific-notation literals are cast via [double]"..." because the COM
# script parser chokes on bare `e-NN` exponent suffixes in numeric literals.
$ws.Range("B2").Value = -0.9956347978830281
$ws.Range("C2").Value = 0.9989819477307562
$ws.Range("D2").Value = 0.9956347978830281
$ws.Range("E2").Value = 0.9993591157448297

$ws.Range("B3").Value = [double]"9.89583079448934e-14"
$ws.Range("C3").Value = [double]"1.603854010522224e-17"
$ws.Range("D3").Value = [double]"9.89583079448934e-14"
$ws.Range("E3").Value = [double]"9.990175297692865e-19"

$ws.Range("D4").Value = -1
$ws.Range("E4").Value = -0.9982011133169525

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = [double]"4.873490752278313e-16"

$ws.Range("B6").Value = -0.9905828645282254
$ws.Range("D6").Value = 0.9905828645282254
$ws.Range("E6").Value = 0.9967268938497555

$ws.Range("B7").Value = [double]"9.867900941422818e-12"
$ws.Range("D7").Value = [double]"9.867900941422818e-12"
$ws.Range("E7").Value = [double]"1.762806016821195e-14"

$ws.Range("D8").Value = 0.9982011133169525

$ws.Range("D9").Value = [double]"4.873490752278313e-16"

# ---- Styling: bold, boxed (thin, all sides) border, centered + top-aligned ----
# Build the combined format once on an off-grid scratch cell, then copy/paste
# only the formatting onto every header/label cell one at a time. Doing this
# per-cell (via PasteSpecial) instead of re-running the four property
# assignments (Font.Bold / Borders.LineStyle / HorizontalAlignment /
# VerticalAlignment) against every destination keeps the saved style table
# to just the two cellXfs actually in use, instead of also persisting every
# transient intermediate combination.
$scratch = $ws.Range("ZZ1000")
$scratch.Value = "x"
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Copy()

$styledCells = @("B1","C1","D1","E1","A2","A3","A4","A5","A6","A7","A8","A9")
foreach ($addr in $styledCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$scratch.Clear()
